function Add-Run {
    param($doc, $pos, $text, $bold, $italic)
    $r = $doc.Range($pos, $pos)
    $r.InsertAfter($text)
    if ($bold) { $r.Font.Bold = 1 }
    if ($italic) { $r.Font.Italic = 1 }
    return $r.End
}

$d = $word.ActiveDocument

# --- Fix dates ---
$d.Content.Find.Execute("6/12/2020", $true, $false, $false, $false, $false, $true, 1, $false, "6/03/2021", 2) | Out-Null
$d.Content.Find.Execute("Version 0.3, 2020-08-13", $true, $false, $false, $false, $false, $true, 1, $false, "Version 0.3, 2021-06-03", 2) | Out-Null

# --- Locate target paragraph (the long "For each of the three variables..." list item) ---
$targetIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text -like "For each of the three variables, in the Graph*") {
        $targetIdx = $i
        break
    }
}

$p = $d.Paragraphs($targetIdx)
$p.Range.Text = ""
$pos = $p.Range.Start
$pos = Add-Run $d $pos 'In the Graph tab in the upper tool bar, make sure the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos '“' $false $false
$pos = Add-Run $d $pos 'Show auxiliary graph' $false $false
$pos = Add-Run $d $pos '”' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'box is checked.' $false $false

# --- Insert the new paragraphs after the target paragraph ---
$insertAfterIdx = $targetIdx
$numNewParagraphs = 9
$d.Paragraphs($insertAfterIdx).Range.InsertParagraphAfter()

# paragraph 1: style=FirstParagraph
$curIdx = $insertAfterIdx + 1
$np = $d.Paragraphs($curIdx)
$np.Style = 'FirstParagraph'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'The primary plot (the large one on the left of the tab) shows a scatter plot of the response variable versus explanatory variable and the best fitting model. Imagine that each of the scattered points were raised or lowered to fall exactly on the best fitting model. This vertical position corresponds to the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'model values' $true $false
$pos = Add-Run $d $pos '.' $false $false

# paragraph 2: style=BodyText
$curIdx = $insertAfterIdx + 2
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'The auxiliary graph is on the right. It shows two clouds of points. The right cloud gives the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'raw values' $true $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'of the response variable. The vertical position of each point in the raw cloud is identical to the vertical position of the corresponding data point in the primary graph.' $false $false

# paragraph 3: style=BodyText
$curIdx = $insertAfterIdx + 3
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'The left cloud in the auxiliary graph gives the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'model values' $true $false
$pos = Add-Run $d $pos '. Again, the vertical position of each point in the model-value cloud is identical to the corresponding model value in the primary plot.' $false $false

# paragraph 4: style=BodyText
$curIdx = $insertAfterIdx + 4
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'Both the raw and model-value clouds in the auxiliary graph are marked with an I-shaped interval. This vertical interval covers approximately 95% of the points in its cloud. The center of the interval is the mean points in the cloud, the ends are plus-or-minus 2 standard deviations away from this.' $false $false

# paragraph 5: style=BodyText
$curIdx = $insertAfterIdx + 5
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'One helpful way to describe a relationship between two variables is to quantify how much of the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'variation' $false $true
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'in the response variable can be accounted for by the explanatory variable. A standard way to quantify this is with a statistic called R-squared, which always falls between 0 and 1. Zero means no relationship and 1 is a' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos '“' $false $false
$pos = Add-Run $d $pos 'perfect' $false $false
$pos = Add-Run $d $pos '”' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'relationship where the explanatory variable exactly accounts for the response variable. Think of R-squared as measuring the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'strength of the relationship' $true $false
$pos = Add-Run $d $pos '. More precisely, R-squared is the fraction of the variance of the response variable accounted for by the explanatory variable.' $false $false

# paragraph 6: style=BodyText
$curIdx = $insertAfterIdx + 6
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'You can estimate the strength of the relationship, R-squared, from the auxiliary graph.' $false $false

# paragraph 7: style=BlockText
$curIdx = $insertAfterIdx + 7
$np = $d.Paragraphs($curIdx)
$np.Style = 'BlockText'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'R is the ratio of the length of the model value interval to the length of the raw values interval. Square R to get R-squared.' $false $false

# paragraph 8: style=FirstParagraph
$curIdx = $insertAfterIdx + 8
$np = $d.Paragraphs($curIdx)
$np.Style = 'FirstParagraph'
$np.Range.InsertParagraphAfter()
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'Another way to describe the relationship between the explanatory and response variables is with the' $false $false
$pos = Add-Run $d $pos ' ' $false $false
$pos = Add-Run $d $pos 'effect size' $true $false
$pos = Add-Run $d $pos '. Whereas R-squared is always on the scale 0 to 1, the effect size reflects the actual units of the explanatory and response variables. It is the change in model value (' $false $false
$pos = Add-Run $d $pos '“' $false $false
$pos = Add-Run $d $pos 'rise' $false $false
$pos = Add-Run $d $pos '”' $false $false
$pos = Add-Run $d $pos ') per unit change in the explanatory variable (' $false $false
$pos = Add-Run $d $pos '“' $false $false
$pos = Add-Run $d $pos 'run' $false $false
$pos = Add-Run $d $pos '”' $false $false
$pos = Add-Run $d $pos '). The ratio is rise over run, in other words, the slope of the model line.' $false $false

# paragraph 9: style=BodyText
$curIdx = $insertAfterIdx + 9
$np = $d.Paragraphs($curIdx)
$np.Style = 'BodyText'
$pos = $np.Range.Start
$pos = Add-Run $d $pos 'For each of the three variables, list the strength of the relationship both as a fraction of the variation explained (R-squared) and as the change in systolic blood pressure per unit change of the explanatory variable (slope of model line).' $false $false

# --- Change the style of the formerly-adjacent "Fill in the table..." paragraph ---
$lastIdx = $insertAfterIdx + 10
$lastP = $d.Paragraphs($lastIdx)
$lastItalicText = "Fill in the table with your answers."
$lastItalicStart = $lastP.Range.Start
$lastP.Style = 'BodyText'
# Changing the paragraph Style strips existing run-level formatting in this runtime;
# restore the italic run that was already there before our edit.
$lastItalicRange = $d.Range($lastItalicStart, $lastItalicStart + $lastItalicText.Length)
$lastItalicRange.Font.Italic = 1

Write-Output "done"
